$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-19 -> 2023-09-20, i.e. 45188 -> 45189) for every data row (2-185).
$ws.Range("C2:C185").Value = 45189
